# Updated cryptos list data refresh (GitHub Actions style price update)
# This script updates the Coin, Link, Price and Volume(1h) columns for the
# crypto listing on Sheet1, rows 2-51, matching a refreshed data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) contains values that look numeric (e.g. "0.6294",
# "1.002") but must be preserved as literal text, exactly as the source
# data feed writes them (leading/trailing zeros, multi-dot big numbers like
# "29.114.95", etc.). Setting NumberFormat to Text ("@") on the whole
# column before assigning the values prevents Excel from silently
# re-interpreting them as numbers. We restore the cell style afterwards so
# no stray formatting is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.114.95"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.834.49"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "243.84"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "0.6294"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "0.07471"
$ws.Range("E8").Value = "  -1.80%  "
$ws.Range("D9").Value = "0.2931"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").Value = "23.06"
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("D11").Value = "0.07726"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("D12").Value = "1.834.21"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "5.003"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").Value = "0.6686"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").Value = "83.09"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "0.000009388"
$ws.Range("E16").Value = "  -4.51%  "
$ws.Range("D17").Value = "6.066"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "29.160.50"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").Value = "12.61"
$ws.Range("E19").Value = "  +2.31%  "
$ws.Range("D20").Value = "224.01"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").Value = "7.131"
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").Value = "160.24"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("D26").Value = "8.513"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").Value = "17.95"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").Value = "1.502"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").Value = "4.149"
$ws.Range("E29").Value = "  +2.11%  "
$ws.Range("D30").Value = "4.066"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").Value = "0.05473"
$ws.Range("E31").Value = "  +5.42%  "
$ws.Range("D32").Value = "1.203"
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").Value = "0.7515"
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("D34").Value = "1.857"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("D36").Value = "2.608"
$ws.Range("E36").Value = "  -3.25%  "
$ws.Range("D37").Value = "1.230.45"
$ws.Range("E37").Value = "  -2.92%  "
$ws.Range("D38").Value = "2.758"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "0.01790"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "6.664"
$ws.Range("E40").Value = "  +6.65%  "
$ws.Range("D41").Value = "0.8953"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "1.004"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").Value = "101.73"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").Value = "65.69"
$ws.Range("E44").Value = "  +1.68%  "
$ws.Range("D45").Value = "0.00000000124"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.5098"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("B47").Value = "XinFinNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D47").Value = "0.07645"
$ws.Range("E47").Value = "  +10.52%  "
$ws.Range("D48").Value = "0.4052"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("D49").Value = "8.978"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.658"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.05804"
$ws.Range("E51").Value = "  +0.93%  "

# Restore the original (default) style now that the values are written, so
# we don't leave an explicit text-format style applied to the cells.
$priceRange.Style = "Normal"
